$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the extent of the data (4 columns x 3 rows in this sheet).
$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

# Capture the current (pre-edit) values of the data block.
$vals = @{}
for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $vals["$r,$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

# Rearrange the columns: the last column moves to the front, and the
# remaining columns shift right by one (new order = old[last], old[1..last-1]).
for ($r = 1; $r -le $rowCount; $r++) {
    $last = $vals["$r,$colCount"]
    $ws.Cells.Item($r, 1).Value = $last
    for ($c = 2; $c -le $colCount; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals["$r,$($c - 1)"]
    }
}

# Update the sheet selection to the new first column's data range.
$ws.Range("A1:A3").Select()
